$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$newRows = @(
    @(44330, 3, 7, 99.5732574679943),
    @(44331, 0, 4, 56.89900426742533),
    @(44332, 0, 3, 42.67425320056899),
    @(44333, 0, 3, 42.67425320056899),
    @(44334, 0, 3, 42.67425320056899),
    @(44335, 1, 4, 56.89900426742533),
    @(44336, 0, 4, 56.89900426742533),
    @(44337, 0, 1, 14.22475106685633),
    @(44338, 0, 1, 14.22475106685633),
    @(44339, 0, 1, 14.22475106685633),
    @(44340, 0, 1, 14.22475106685633),
    @(44341, 0, 1, 14.22475106685633),
    @(44342, 0, 0, 0),
    @(44343, 0, 0, 0)
)

$startRow = 256
$lastFormattedRow = 255

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    # Copy the date column's style (border/font/alignment/number format) from the
    # last existing row so the new A-column cells stay visually consistent.
    $ws.Range("A$lastFormattedRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
}

$excel.CutCopyMode = 0
